$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Cells.Item(8, 1).Value = "Volume 32   Number  24"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Data table updates (rows 14-31) ---
$ws.Range("N14").Value = -88.888888888888
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -36.363636363636
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 92
$ws.Range("K16").Value = -10.869565217391
$ws.Range("L16").Value = -1.204819277108
$ws.Range("M16").Value = -44.594594594594
$ws.Range("N16").Value = -84.926470588235
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -51.351351351351
$ws.Range("I17").Value = 124
$ws.Range("J17").Value = 141
$ws.Range("K17").Value = -12.056737588652
$ws.Range("L17").Value = -9.48905109489
$ws.Range("M17").Value = 49.397590361445
$ws.Range("N17").Value = -45.614035087719
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -60
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -36
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 164
$ws.Range("K18").Value = -31.70731707317
$ws.Range("L18").Value = 17.894736842105
$ws.Range("M18").Value = -34.117647058823
$ws.Range("N18").Value = -80.212014134275
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -29.230769230769
$ws.Range("I19").Value = 311
$ws.Range("J19").Value = 331
$ws.Range("K19").Value = -6.042296072507
$ws.Range("L19").Value = -5.182926829268
$ws.Range("M19").Value = 66.310160427807
$ws.Range("N19").Value = 34.051724137931
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 7.142857142857
$ws.Range("I20").Value = 61
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 1.666666666666
$ws.Range("L20").Value = -12.857142857142
$ws.Range("M20").Value = -11.59420289855
$ws.Range("N20").Value = -83.819628647214
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -28.947368421052
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 166
$ws.Range("H21").Value = -33.132530120481
$ws.Range("I21").Value = 698
$ws.Range("J21").Value = 802
$ws.Range("K21").Value = -12.967581047381
$ws.Range("L21").Value = -3.591160220994
$ws.Range("M21").Value = 5.597579425113
$ws.Range("N21").Value = -64.658227848101
$ws.Range("C22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 6.666666666666
$ws.Range("M22").Value = -20
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -7.142857142857
$ws.Range("I23").Value = 70
$ws.Range("J23").Value = 85
$ws.Range("K23").Value = -17.647058823529
$ws.Range("L23").Value = -21.348314606741
$ws.Range("M23").Value = 16.666666666666
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -14.912280701754
$ws.Range("I24").Value = 514
$ws.Range("J24").Value = 503
$ws.Range("K24").Value = 2.186878727634
$ws.Range("L24").Value = 7.531380753138
$ws.Range("M24").Value = -2.281368821292
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 25
$ws.Range("H25").Value = 13.636363636363
$ws.Range("I25").Value = 145
$ws.Range("J25").Value = 171
$ws.Range("K25").Value = -15.204678362573
$ws.Range("L25").Value = 150
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = -40.74074074074
$ws.Range("I26").Value = 171
$ws.Range("J26").Value = 266
$ws.Range("K26").Value = -35.714285714285
$ws.Range("L26").Value = -17.391304347826
$ws.Range("M26").Value = -21.917808219178
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = -50
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = 8
$ws.Range("L28").Value = -10
$ws.Range("N29").Value = -93.75
$ws.Range("N30").Value = -93.548387096774
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("G31").Value = 2
$ws.Range("L31").Value = 140
